$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: add "alleen" in column A (existing Testknop row) ---
$ws.Range("A4").Value = "alleen"

# --- Row 5: add "alleen" in column A (existing criculair buffer row) ---
$ws.Range("A5").Value = "alleen"
# normalize D5:E5 number format (was a duplicate HH:MM format, now shares the
# canonical HH:MM format used elsewhere in the sheet)
$ws.Range("D5:E5").NumberFormat = "HH:MM"

# --- Row 6: new log entry ---
$ws.Range("A6").Value = "alleen"
$ws.Range("B6").Value = "parser werkend maar nog niet geïnplementeerd"
$ws.Range("C6").Value = 43385
$ws.Range("C6").NumberFormat = "DD/MM/YY"
$ws.Range("D6").Value = 0.375
$ws.Range("E6").Value = 0.5
$ws.Range("D6:E6").NumberFormat = "HH:MM:SS"

# --- Row 7: new log entry ---
$ws.Range("A7").Value = "alleen"
$ws.Range("B7").Value = "parser werkend en geïnplementeerd, maar alleen nog de basis foutmeldingen"
$ws.Range("C7").Value = 43396
$ws.Range("C7").NumberFormat = "DD/MM/YY"
$ws.Range("D7").Value = 0.53125
$ws.Range("E7").Value = 0.625
$ws.Range("D7:E7").NumberFormat = "HH:MM:SS"

# --- widen column B to fit the longer log text ---
$ws.Columns.Item(2).ColumnWidth = 65.75

# --- restore the active selection like in the authored workbook ---
$ws.Range("F12").Select()
